# Merge the three split runs "<id>" + "p114r_1" + "</id>" (inside the
# paragraph that currently reads "<id>p114r_1</id>") back into a single
# run, keeping the formatting of the first ("<id>") run - Courier New,
# color 7f6000, sz/szCs 18, rtl 0.
#
# A formatted Find & Replace over the whole document content does this:
# Word's Find/Replace merges any runs spanned by a match into one run
# that carries the formatting of the first run in the match, which is
# exactly the target state shown by the diff. The search text is unique
# in the document (other "<id>...</id>" tags wrap different ids, e.g.
# "fig_p114r_1"), so this only touches the intended paragraph.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "<id>p114r_1</id>",  # Find What
    $false,              # MatchCase
    $true,               # MatchWholeWord
    $false,              # MatchWildcards
    $false,              # MatchSoundsLike
    $false,              # MatchAllWordForms
    $true,               # Forward
    1,                   # Wrap (wdFindContinue)
    $false,              # Format
    "<id>p114r_1</id>",  # Replace With
    2                    # Replace (wdReplaceAll)
)
